$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as Text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.808.98'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '2.290.65'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '116.55'
$ws.Range("E5").Value = '  +1.90%  '

$ws.Range("D6").Value = '267.36'
$ws.Range("E6").Value = '  -1.31%  '

$ws.Range("E7").Value = '  +2.62%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  -1.00%  '

$ws.Range("D10").Value = '48.04'
$ws.Range("E10").Value = '  -0.41%  '

$ws.Range("D11").Value = '0.0940'
$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("D12").Value = '9.18'
$ws.Range("E12").Value = '  +0.96%  '

$ws.Range("E13").Value = '  +1.28%  '

$ws.Range("D14").Value = '15.55'
$ws.Range("E14").Value = '  -2.22%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.893'
$ws.Range("E15").Value = '  +4.39%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.633.34'
$ws.Range("E16").Value = '  -0.55%  '

$ws.Range("D17").Value = '2.289.46'
$ws.Range("E17").Value = '  -0.55%  '

$ws.Range("D18").Value = '43.694.92'
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("E19").Value = '  -0.38%  '

$ws.Range("D20").Value = '6.95'
$ws.Range("E20").Value = '  +1.78%  '

$ws.Range("D21").Value = '72.50'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '2.47'
$ws.Range("E22").Value = '  +0.44%  '

$ws.Range("D23").Value = '236.13'
$ws.Range("E23").Value = '  +1.16%  '

$ws.Range("D24").Value = '9.67'
$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").Value = '2.89'
$ws.Range("E25").Value = '  -2.26%  '

$ws.Range("E26").Value = '  +1.76%  '

$ws.Range("D27").Value = '11.78'
$ws.Range("E27").Value = '  +0.30%  '

$ws.Range("D28").Value = '42.28'
$ws.Range("E28").Value = '  +0.80%  '

$ws.Range("D29").Value = '3.42'
$ws.Range("E29").Value = '  +0.65%  '

$ws.Range("E30").Value = '  -0.56%  '

$ws.Range("D31").Value = '174.06'
$ws.Range("E31").Value = '  -0.86%  '

$ws.Range("D32").Value = '21.80'
$ws.Range("E32").Value = '  +0.73%  '

$ws.Range("D33").Value = '0.0916'
$ws.Range("E33").Value = '  -2.30%  '

$ws.Range("D34").Value = '5.75'
$ws.Range("E34").Value = '  +0.50%  '

$ws.Range("D35").Value = '0.131'
$ws.Range("E35").Value = '  +2.15%  '

$ws.Range("E36").Value = '  +5.27%  '

$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("E38").Value = '  +3.31%  '

$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("D40").Value = '2.57'
$ws.Range("E40").Value = '  +7.80%  '

$ws.Range("D41").Value = '14.16'
$ws.Range("E41").Value = '  +3.14%  '

$ws.Range("D42").Value = '74.40'
$ws.Range("E42").Value = '  +0.04%  '

$ws.Range("E43").Value = '  -2.88%  '

$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = '5.98'
$ws.Range("E44").Value = '  -7.06%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.21%  '

$ws.Range("E47").Value = '  +3.57%  '

$ws.Range("D48").Value = '8.64'
$ws.Range("E48").Value = '  -2.76%  '

$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("D50").Value = '72.25'
$ws.Range("E50").Value = '  +32.93%  '

$ws.Range("D51").Value = '101.79'
$ws.Range("E51").Value = '  -0.92%  '
